$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A595").Value = "TEST"
